# Fix Training Data Issue (#48)
# The "Date" column (BF) held a malformed value ("6-17-2007-08") for every
# data row. NBA.com's stat-API season/date display caused the date to be
# copied in a mangled "M-D-YYYY-YY" form instead of the real game date. The
# underlying game actually took place on 2008-06-17, so every row's BF cell
# needs to be corrected to that (plain, textual) date value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "6-17-2007-08"
$newValue = "2008-06-17"

# Column BF ("Date") runs from row 2 through row 31.
$firstRow = 2
$lastRow  = 31
$bfCol    = 58   # column BF

$rng = $ws.Range($ws.Cells.Item($firstRow, $bfCol), $ws.Cells.Item($lastRow, $bfCol))

# Force the cells to be treated as plain text so the replacement keeps the
# literal "2008-06-17" string instead of Excel auto-converting it to a date
# serial number.
$rng.NumberFormat = "@"

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $bfCol)
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value = $newValue
    }
}

# Drop the temporary text format again so the cells keep their original
# (default) style, matching the rest of the sheet.
$rng.Style = "Normal"
